# Update "想去人数" (interest count) values in the "展览" and "全部类型" sheets
# to match the freshly generated data for commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 14151
$ws1.Range("F3").Value = 558
$ws1.Range("F5").Value = 1219
$ws1.Range("F7").Value = 13950
$ws1.Range("F8").Value = 15071
$ws1.Range("F10").Value = 21
$ws1.Range("F19").Value = 65
$ws1.Range("F21").Value = 1169
$ws1.Range("F22").Value = 124
$ws1.Range("F24").Value = 5857
$ws1.Range("F26").Value = 1071
$ws1.Range("F27").Value = 5470
$ws1.Range("F28").Value = 61
$ws1.Range("F29").Value = 132
$ws1.Range("F30").Value = 71
$ws1.Range("F31").Value = 358

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 14151
$ws4.Range("F4").Value = 558
$ws4.Range("F6").Value = 1219
$ws4.Range("F8").Value = 13950
$ws4.Range("F9").Value = 15072
$ws4.Range("F11").Value = 21
$ws4.Range("F20").Value = 65
$ws4.Range("F22").Value = 1169
$ws4.Range("F23").Value = 124
$ws4.Range("F26").Value = 5857
$ws4.Range("F28").Value = 1071
$ws4.Range("F29").Value = 5470
$ws4.Range("F30").Value = 61
$ws4.Range("F31").Value = 132
$ws4.Range("F32").Value = 71
$ws4.Range("F33").Value = 358
